$d = $word.ActiveDocument

# Update the date heading
$d.Paragraphs.Item(1).Range.Find.Execute("2022-12-18 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2022-12-19 Monday", 2) | Out-Null

# Update each arithmetic expression cell in the table, in reading order
$t = $d.Tables.Item(1)
$values = @(
    "16+54=",
    "21+67=",
    "15+45=",
    "74-17=",
    "28-7=",
    "34+28=",
    "3+34=",
    "33+62=",
    "50+17=",
    "50+22=",
    "15+31=",
    "9+26=",
    "41-24=",
    "10-2=",
    "80-56=",
    "4+51=",
    "30+20=",
    "17+21=",
    "95-67=",
    "61-24=",
    "31+45=",
    "72-13=",
    "88-54=",
    "97-60=",
    "32+67=",
    "46-8=",
    "20-19=",
    "80-46=",
    "10+0=",
    "97-11=",
    "50-46=",
    "13+7=",
    "36+47=",
    "94-10=",
    "93+0=",
    "92-86=",
    "77-30=",
    "40+58=",
    "15+5=",
    "8+71=",
    "33-15=",
    "17-7=",
    "88+10=",
    "12+71=",
    "81-40=",
    "8+91=",
    "4+70=",
    "25-9=",
    "2+26=",
    "43-34=",
    "64+5=",
    "37+46=",
    "38-7=",
    "76-58=",
    "23+0=",
    "67+3=",
    "36-14=",
    "93-0=",
    "53+13=",
    "25+30=",
    "50-7=",
    "94-63=",
    "65+2=",
    "61+7=",
    "76-0=",
    "81-50=",
    "36+28=",
    "77-4=",
    "78+6=",
    "47+46=",
    "69+30=",
    "68+1=",
    "21-9=",
    "40-4=",
    "35+62=",
    "65-17=",
    "69+30=",
    "22+6=",
    "65-30=",
    "34+0=",
    "4+70=",
    "33+52=",
    "40+6=",
    "95-62=",
    "91-37=",
    "95-36=",
    "1+7=",
    "45-6=",
    "49-39=",
    "51+22=",
    "10+59=",
    "52-6=",
    "82-63=",
    "0+99=",
    "42-17=",
    "81-66=",
    "92-36=",
    "63+0=",
    "74+5=",
    "2+61="
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cellRange = $cell.Range
        $cellRange.MoveEnd(1, -1) | Out-Null
        $cellRange.Text = $values[$idx]
        $idx = $idx + 1
    }
}
